# Add new parameter rows (10-16) to the Uncaging result-table template,
# describing additional fields used for uncaging experiments.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: laser_power / double / laser power (percent)
$ws.Range("A10").Value = "laser_power"
$ws.Range("B10").Value = "double"
$ws.Range("C10").Value = "laser power (percent)"

# Row 11: laser_wavelength / double / laser wavelength
$ws.Range("A11").Value = "laser_wavelength"
$ws.Range("B11").Value = "double"
$ws.Range("C11").Value = "laser wavelength"

# Row 12: shutter_open / string / T or F for shutter open
$ws.Range("A12").Value = "shutter_open"
$ws.Range("B12").Value = "string"
$ws.Range("C12").Value = "T or F for shutter open"

# Row 13: drug_condition / string / drugs in the bath (free text)
$ws.Range("A13").Value = "drug_condition"
$ws.Range("B13").Value = "string"
$ws.Range("C13").Value = "drugs in the bath (free text)"

# Rows 14-15: field names first (number_of_sequences, number_of_stim_groups)
$ws.Range("A14").Value = "number_of_sequences"
$ws.Range("B14").Value = "uint16"
$ws.Range("A15").Value = "number_of_stim_groups"
$ws.Range("B15").Value = "uint16"

# Rows 14-15: descriptions filled in afterwards
$ws.Range("C14").Value = "number of repeats "
$ws.Range("C15").Value = "number of different uncaging locations"

# Row 16: group_names / string / free text used to associate ROI files with this epoch
$ws.Range("A16").Value = "group_names"
$ws.Range("B16").Value = "string"
$ws.Range("C16").Value = "free text used to associate ROI files with this epoch"

# Move the view to match the edited selection/scroll position
# (selection -> A16, scrolled so row 8 is the top visible row).
$ws.Range("A16").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1
